{"js": "// Rename \"Opensense\"/\"OpenSense\" -> \"Opnsense\"/\"OpnSense\" throughout the\n// body (the 4 occurrences the diff touches), and relocate the \"_GoBack\"\n// bookmark from the \"(See Figure 4)\" paragraph to the 3rd-from-last\n// (empty) paragraph at the end of the document, right before the final\n// sectPr.\n\nconst body = context.document.body;\n\n// 1) Fix the lower-case-\"s\" spelling: \"Opensense\" -> \"Opnsense\".\nconst lower = body.search(\"Opensense\", { matchCase: true, matchWholeWord: false });\nlower.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < lower.items.length; i++) {\n  lower.items[i].insertText(\"Opnsense\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) Fix the capital-\"S\" spelling: \"OpenSense\" -> \"OpnSense\".\nconst upper = body.search(\"OpenSense\", { matchCase: true, matchWholeWord: false });\nupper.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < upper.items.length; i++) {\n  upper.items[i].insertText(\"OpnSense\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 3) Move the \"_GoBack\" bookmark down to the 3rd-from-last paragraph.\ncontext.document.deleteBookmark(\"_GoBack\");\n\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst target = paragraphs.items[paragraphs.items.length - 3];\ntarget.getRange().insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Rename \"Opensense\"/\"OpenSense\" -> \"Opnsense\"/\"OpnSense\" throughout the\n# body (the 4 occurrences the diff touches), and relocate the \"_GoBack\"\n# bookmark from the \"(See Figure 4)\" paragraph to the 3rd-from-last\n# (empty) paragraph at the end of the document, right before the final\n# sectPr.\n\n$d = $word.ActiveDocument\n\n# 1) Fix the lower-case-\"s\" spelling: \"Opensense\" -> \"Opnsense\".\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Text = \"Opensense\"\n$find1.MatchCase = $true\n$find1.MatchWholeWord = $false\n$find1.Replacement.ClearFormatting()\n$find1.Replacement.Text = \"Opnsense\"\n$find1.Execute($null, $true, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n# 2) Fix the capital-\"S\" spelling: \"OpenSense\" -> \"OpnSense\".\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"OpenSense\"\n$find2.MatchCase = $true\n$find2.MatchWholeWord = $false\n$find2.Replacement.ClearFormatting()\n$find2.Replacement.Text = \"OpnSense\"\n$find2.Execute($null, $true, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n# 3) Move the \"_GoBack\" bookmark down to the 3rd-from-last paragraph.\n$paraCount = $d.Paragraphs.Count\n$target = $d.Paragraphs.Item($paraCount - 2).Range\n$d.Bookmarks.Add(\"_GoBack\", $target) | Out-Null\n"}
